# Updated cryptos list data (prices / 1h volume %) per GitHub Actions scrape refresh.
# Price cells (column D) that look like plain numbers are forced to Text format
# first so Excel's COM layer doesn't silently coerce them into numeric values
# (which would drop significant trailing/leading zeros, e.g. "91.40" -> 91.4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.868.02"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "1.875.45"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.45"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4832"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07368"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9397"
$ws.Range("E10").Value = "  +0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.99"
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07780"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "1.894.03"
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.525"
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.592"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.40"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008872"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "27.875.10"
$ws.Range("E20").Value = "  +2.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.84"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.121"
$ws.Range("D23").Value = "2.120.68"
$ws.Range("E23").Value = "  +2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.91"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.52"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.56"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.042"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.96"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.969"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08885"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.336"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7725"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.648"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.731"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02045"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5616"
$ws.Range("E39").Value = "  +5.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05376"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.041"
$ws.Range("E42").Value = "  -0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.536"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4881"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.65"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.56"
$ws.Range("E47").Value = "  +2.91%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.667"
$ws.Range("E49").Value = "  +2.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.15"
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("E51").Value = "  +0.68%  "
